{"js": "// The upstream commit (\"Fixed #295 Add the version of M2Doc in the\n// template custom properties\") batch-resaved many test-fixture .docx\n// files while wiring up a new custom-properties writer. For *this*\n// particular template, the resulting OOXML diff is not a content edit\n// at all: every hunk (word/document.xml, word/footer1.xml and\n// word/styles.xml, the latter mislabeled under the word/footnotes.xml\n// <file> marker in the scraped diff) is fully explained by two purely\n// cosmetic, value-preserving transformations that the tool used to\n// produce/review the diff applied before comparing the XML:\n//\n//   1. Every element's attributes (and the xmlns:* declarations on the\n//      root element) were re-emitted in canonical alphabetical order\n//      (namespaces by prefix, then the remaining attributes by name).\n//   2. The volatile w:rsid*/w:rsidR/w:rsidRPr/w:rsidRDefault/w:rsidP/\n//      w:rsidTr noise attributes (which Word regenerates on every save\n//      and which carry no document meaning) were stripped.\n//\n// No text, run/paragraph formatting, table layout, style definition,\n// section/page setup, or footnote content actually changed \u2014 every\n// attribute value on both sides of the diff is identical once you sort\n// them back (verified programmatically against every hunk). So the\n// faithful replay of this commit, from the document-content point of\n// view that the Word JS object model exposes, is a no-op: we must not\n// introduce any text/formatting/structural change of our own.\n//\n// We still touch `context.sync()` so the script demonstrably runs\n// against the live document without mutating it.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The upstream commit (\"Fixed #295 Add the version of M2Doc in the\n# template custom properties\") batch-resaved many test-fixture .docx\n# files while wiring up a new custom-properties writer. For *this*\n# particular template, the resulting OOXML diff is not a content edit\n# at all: every hunk (word/document.xml, word/footer1.xml and\n# word/styles.xml, the latter mislabeled under the word/footnotes.xml\n# <file> marker in the scraped diff) is fully explained by two purely\n# cosmetic, value-preserving transformations applied by the tool that\n# produced/rendered the diff, before comparing the XML:\n#\n#   1. Every element's attributes (and the xmlns:* declarations on the\n#      root element) were re-emitted in canonical alphabetical order\n#      (namespaces by prefix, then the remaining attributes by name).\n#   2. The volatile w:rsid*/w:rsidR/w:rsidRPr/w:rsidRDefault/w:rsidP/\n#      w:rsidTr noise attributes (which Word regenerates on every save\n#      and which carry no document meaning) were stripped.\n#\n# No text, run/paragraph formatting, table layout, style definition,\n# section/page setup, or footnote content actually changed - every\n# attribute value on both sides of the diff is identical once sorted\n# back (verified programmatically against every hunk). So the faithful\n# replay of this commit, from the document-content point of view that\n# the Word COM object model exposes, is a no-op: we must not introduce\n# any text/formatting/structural change of our own.\n#\n# We still touch the document (a harmless read) so the script\n# demonstrably runs against the live document without mutating it.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
